# LOM3094.xlsx update
# - Fills in missing Portuguese body text for Objetivos, Programa resumido,
#   Programa, Metodo, Criterio, Norma de recuperacao and Bibliografia.
# - Inserts a new row for "Docentes responsaveis:" teacher name (previously
#   misplaced inside the "Objetivos"/"Metodo" rows), shifting every row
#   below it down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Objetivos (row 10): replace the misplaced teacher name with the real
#    Portuguese objectives text.
$objetivos = "1. Compreensão dos fenômenos envolvidos no processo de solidificação. 2. Noção dos aspectos práticos do processamento de materiais em empresas de fundição."
$ws.Range("B10").Value = $objetivos
$ws.Range("C10").Value = $objetivos

# 2) Insert a new row at 13 (pushes old rows 13-23 down to 14-24) to hold the
#    teacher name under "Docentes responsáveis:" (row 12).
$ws.Rows.Item(13).Insert()

$professor = "5009972 - Gilberto Carvalho Coelho"
$ws.Range("B13").Value = $professor
$ws.Range("C13").Value = $professor
$ws.Range("B13").Style = $ws.Range("B19").Style
$ws.Range("C13").Style = $ws.Range("C19").Style

# 3) Programa resumido (now row 14): replace placeholder "Semestral" with the
#    real short-syllabus text.
$resumido = "A. Introdução; B. Nucleação; C. Redistribuição de soluto na solidificação de ligas; D. Crescimento; E. Macroestruturas de solidificação; F. Aspectos práticos do processamento de materiais por fundição."
$ws.Range("B14").Value = $resumido
$ws.Range("C14").Value = $resumido

# 4) Programa (now row 16): replace placeholder date with the full syllabus
#    text.
$programa = "1. Introdução: história da fundição; aplicações e mercado de fundidos; 2. Nucleação: Nucleação homogênea; nucleação heterogênea; taxa de nucleação, agentes nucleantes; 3. Redistribuição de soluto na solidificação: Materiais puros; ligas binárias; coeficiente de redistribuição; solidificação em condições de equilíbrio; solidificação fora de equilíbrio; 4. Crescimento: solidificação de ligas monofásicas - crescimento planar, celular e dendrítico; solidificação de ligas polifásicas ligas eutéticas e ligas peritéticas; 5. Macroestruturas de solidificação: contração volumétrica na solidificação; zonas coquilhada, colunar e equiaxial; controle da macroestrutura; 6. Aspectos práticos do processamento de materiais por fundição: equipamentos e processos de fundição; segregação macro e microssegregação; defeitos originados na solidificação."
$ws.Range("B16").Value = $programa
$ws.Range("C16").Value = $programa

# 5) Metodo (now row 19): replace misplaced teacher name with the real
#    method text.
$metodo = "O curso será ministrado na forma de aulas expositivas. Estão previstas visitas a empresas de fundição para consolidação dos conceitos apresentados nas aulas expositivas."
$ws.Range("B19").Value = $metodo
$ws.Range("C19").Value = $metodo

# 6) Criterio (now row 20): the grading criteria text that used to sit one
#    row up.
$criterio = "Serão aplicadas duas avaliações escritas (P1 e P2) que comporão a nota final (NF). O critério para a nota final é: NF=(P1+P2)/2"
$ws.Range("B20").Value = $criterio
$ws.Range("C20").Value = $criterio

# 7) Norma de recuperacao (now row 21): the make-up exam text that used to
#    sit one row up.
$norma = "Será aplicada uma prova de recuperação cuja nota comporá média aritmética com a nota final NF."
$ws.Range("B21").Value = $norma
$ws.Range("C21").Value = $norma

# 8) Bibliografia (now row 22): fill in the actual bibliography text.
$biblio = "1. Garcia, A. Solidificação: Fundamentos e Aplicações, Editora da Unicamp, 2001. 2. Flemings, M.C. Solidification Processing, McGraw-Hill, 1974. 3. Pfann, W.G. Zone Melting, John Wiley, 1966. 4. Shewmon, P.G. Diffusion in Solids, McGraw-Hill, 1963. 5. Shewmon, P.G. Transformations in Metals, McGraw-Hill, 1969. 6. Prates, M.; Davis, G.J. Solidificação e Fundição de Metais e suas Ligas, EDUSP, 1978. 7. Davis, G.J. Solidification and Casting, Applied Science Publisher, 1973. 8. Brice, J.C. The Growth of Crystals from the Melt, John Wiley, 1965. 9. Winegard, W.C. An Introduction to Solidification of Metals, Institute of Metals, 1964. 10. Chalmers, B. Principles of Solidification, Robert E. Krieger, 1964. 11. Casting, ASM Handbook, Vol 15, Ninth Edition, ASM International, 1988. 12. Metallography and Microstructures, ASM Handbook, Vol 9, Ninth Edition, ASM International, 1988. 13. Welding, Brasing, and Soldering, ASM Handbook, Vol 6, Ninth Edition, ASM International, 1988."
$ws.Range("B22").Value = $biblio
$ws.Range("C22").Value = $biblio

# 9) Cosmetic cleanup of the column definitions: column A only needs its own
#    width entry (column B already carries its own explicit width/style).
$ws.Columns.Item(2).ColumnWidth = $ws.Columns.Item(2).ColumnWidth
